$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 10.862361689314906
$ws.Range("C2").Value = 10.439420823763685
$ws.Range("D2").Value = 8.395051113995466
$ws.Range("E2").Value = 0.37687517840855428

$ws.Range("B3").Value = 28.298408764636484
$ws.Range("C3").Value = 3.2310767151311972
$ws.Range("D3").Value = 1.7686715961889599
$ws.Range("E3").Value = 2.3353231463272066

$ws.Range("B1:E3").Select()
